$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update conversation names in column A
$ws.Range("A2").Value = "conversation_11_07_2023__14_51_17"
$ws.Range("A3").Value = "conversation_11_07_2023__13_22_32"
$ws.Range("A4").Value = "conversation_11_07_2023__13_09_34"

# Clear the now-unused data columns (B:E) for rows 2-4
$ws.Range("B2:E4").ClearContents()
